$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.569.02"
$ws.Range("E2").Value = "  -0.71%  "
$ws.Range("D3").Value = "'2.071.68"
$ws.Range("E3").Value = "  -0.61%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'231.56"
$ws.Range("E5").Value = "  -0.87%  "
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'58.17"
$ws.Range("E8").Value = "  -1.95%  "
$ws.Range("E9").Value = "  -1.92%  "
$ws.Range("D10").Value = "'0.0775"
$ws.Range("E10").Value = "  -1.66%  "
$ws.Range("E11").Value = "  +1.52%  "
$ws.Range("B12").Value = "Chainlink"
$ws.Range("C12").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D12").Value = "'14.84"
$ws.Range("E12").Value = "  +0.52%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "'2.379.70"
$ws.Range("E13").Value = "  -0.52%  "
$ws.Range("D14").Value = "'21.29"
$ws.Range("E14").Value = "  +0.12%  "
$ws.Range("D15").Value = "'0.765"
$ws.Range("E15").Value = "  -1.11%  "
$ws.Range("D16").Value = "'5.33"
$ws.Range("E16").Value = "  -0.29%  "
$ws.Range("D17").Value = "'2.072.54"
$ws.Range("E17").Value = "  -2.56%  "
$ws.Range("D18").Value = "'37.559.14"
$ws.Range("E18").Value = "  -0.59%  "
$ws.Range("D19").Value = "'6.18"
$ws.Range("E19").Value = "  +0.51%  "
$ws.Range("D20").Value = "'70.01"
$ws.Range("E20").Value = "  -2.22%  "
$ws.Range("D21").Value = "'0.0₃0826"
$ws.Range("E21").Value = "  -2.86%  "
$ws.Range("D22").Value = "'227.31"
$ws.Range("E22").Value = "  -0.30%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("E24").Value = "  +0.33%  "
$ws.Range("E25").Value = "  -2.59%  "
$ws.Range("D26").Value = "'9.91"
$ws.Range("E26").Value = "  +6.94%  "
$ws.Range("D27").Value = "'169.64"
$ws.Range("E27").Value = "  -1.18%  "
$ws.Range("D28").Value = "'0.131"
$ws.Range("E28").Value = "  -4.51%  "
$ws.Range("D29").Value = "'19.32"
$ws.Range("E29").Value = "  -0.91%  "
$ws.Range("D31").Value = "'0.121"
$ws.Range("E31").Value = "  +0.51%  "
$ws.Range("E32").Value = "  -3.53%  "
$ws.Range("D33").Value = "'0.0626"
$ws.Range("E33").Value = "  -1.09%  "
$ws.Range("D34").Value = "'4.66"
$ws.Range("E34").Value = "  -1.01%  "
$ws.Range("D35").Value = "'2.54"
$ws.Range("E35").Value = "  +1.21%  "
$ws.Range("E36").Value = "  +0.24%  "
$ws.Range("D37").Value = "'3.30"
$ws.Range("E37").Value = "  -4.00%  "
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("D39").Value = "'5.31"
$ws.Range("E39").Value = "  -1.83%  "
$ws.Range("E40").Value = "  +3.42%  "
$ws.Range("D41").Value = "'98.45"
$ws.Range("E41").Value = "  -1.03%  "
$ws.Range("E42").Value = "  -2.51%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "'1.485.84"
$ws.Range("E43").Value = "  +2.44%  "
$ws.Range("D44").Value = "'2.91"
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").Value = "'1.20"
$ws.Range("E45").Value = "  +3.70%  "
$ws.Range("D46").Value = "'16.71"
$ws.Range("E46").Value = "  -2.43%  "
$ws.Range("B47").Value = "FTXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D47").Value = "'4.03"
$ws.Range("E47").Value = "  -3.11%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "'1.04"
$ws.Range("E48").Value = "  -2.49%  "
$ws.Range("E49").Value = "  -1.64%  "
$ws.Range("D50").Value = "'2.95"
$ws.Range("E50").Value = "  -1.55%  "
$ws.Range("D51").Value = "'2.263.57"
$ws.Range("E51").Value = "  -0.61%  "
